$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the "取得日時" (retrieved at) timestamp for the new scrape batch
# (rows 2-7) from 18:24:13 to 18:33:12 on 2025-10-09.
$ws.Range("A2:A7").Value = "2025-10-09 18:33:12"
